# "Generate Report for Handback" — localization-status.xlsx update
#
# Summary of the change being applied (see xml diff):
#   - Status text everywhere changes from "Ready for handoff" to
#     "Handed back: in sync with en-US" (Overview!E/F and the per-locale
#     "Status" column on the zh-cn / de-de sheets all share this string).
#   - zh-cn sheet: "Latest Target File" (I) / "Latest Handback File" (J) /
#     "Latest Handback DateTime" (K) get populated for rows 2 and 3, and a
#     new hyperlink ("a.md") is added on column I.
#   - de-de sheet: same shape of change, with de-de specific handback file
#     name and a later handback timestamp.
#   - A handful of column widths widen (to fit the new, longer cell text).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1) Status column text: "Ready for handoff" -> "Handed back: in sync
#    with en-US" everywhere it appears (Overview E/F, zh-cn C, de-de C).
# ---------------------------------------------------------------------
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2) zh-cn: Latest Target File (I) / Latest Handback File (J) /
#    Latest Handback DateTime (K) for rows 2 & 3, plus the new
#    hyperlink on column I (pointing at the same a.md source doc).
# ---------------------------------------------------------------------
$aMdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f410ceb6e1d36ec87a166dc3987e15bb7eb0daf0/e2e/a.md"

$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $aMdUrl, $null, $null, "a.md")
$zhcn.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-20 08:45:53"

$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $aMdUrl, $null, $null, "a.md")
$zhcn.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-20 08:45:53"

# ---------------------------------------------------------------------
# 3) de-de: same shape of change, de-de handback file + later timestamp.
# ---------------------------------------------------------------------
$dede.Hyperlinks.Add($dede.Range("I2"), $aMdUrl, $null, $null, "a.md")
$dede.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K2").Value = "2016-08-20 08:45:59"

$dede.Hyperlinks.Add($dede.Range("I3"), $aMdUrl, $null, $null, "a.md")
$dede.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$dede.Range("K3").Value = "2016-08-20 08:45:59"

# ---------------------------------------------------------------------
# 4) Column widths widen to fit the newly-populated / longer text.
#    Excel's ColumnWidth setter only lands on 1/6-character pixel
#    boundaries, so these inputs are chosen to round to the nearest
#    achievable stored width (29.98 -> 30, 40 -> 40 exactly).
# ---------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth  = 29.14   # Overview!E  -> ~29.98
$overview.Columns.Item(6).ColumnWidth  = 29.14   # Overview!F  -> ~29.98

$zhcn.Columns.Item(3).ColumnWidth  = 29.14        # zh-cn!C (Status) -> ~29.98
$zhcn.Columns.Item(10).ColumnWidth = 39.14        # zh-cn!J (Handback File) -> 40

$dede.Columns.Item(3).ColumnWidth  = 29.14        # de-de!C (Status) -> ~29.98
$dede.Columns.Item(10).ColumnWidth = 39.14        # de-de!J (Handback File) -> 40

Write-Output "Generate Report for Handback: applied."
